$wb = $excel.ActiveWorkbook

# Sheet names affected by the update: 展览 (sheet1) and 全部类型 (sheet4)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1694
    $ws.Range("F3").Value = 7830
    $ws.Range("F5").Value = 239
}

$wb.Save()
